# Adds a new "Speed" worksheet computing wheel speed from axle data,
# and makes it the active sheet (mirrors manual edit + "added ADC conversion, added logs").

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and rename it "Speed".
$speed = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$speed.Name = "Speed"

# Header / labels
$speed.Range("B2").Value = "with the G40plus axle"

$speed.Range("C3").Value = "wheel diameter [cm]"
$speed.Range("E3").Value = 40

$speed.Range("C4").Value = "circumference [m]"
$speed.Range("E4").Formula = "=E3/100*3.14"

$speed.Range("C6").Value = "turns at full speed [/s]"
$speed.Range("E6").Value = 2

$speed.Range("C8").Value = "speed [m/s]"
$speed.Range("E8").Formula = "=E4*E6"

$speed.Range("C9").Value = "speed [km/h]"
$speed.Range("E9").Formula = "=E8*3.6"

$speed.Range("E7").Select() | Out-Null

$wb.Worksheets.Item("Speed").Activate() | Out-Null
